$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Modules section, point 1 example: add "at least" before each "one for ..."
#    (commit message: 'Modules section point 1 - added "at least"')
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "One for temperature sensor, one for RTC and one for UART communication.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "At least one for temperature sensor, at least one for RTC and at least one for UART communication.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Add a red note paragraph right after the RTC variable naming example,
#    pointing the reader to the naming conventions for the global/local suffix.
# ---------------------------------------------------------------------------
$rtcRange = $d.Content
$rtcFound = $rtcRange.Find.Execute(
    "RTC_<variable unit identifier ><variable description>_<global or local suffix>",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rtcFound) {
    $rtcPara = $rtcRange.Paragraphs(1)
    $rtcPara.Range.InsertParagraphAfter()
    $notePara = $rtcPara.Next()
    $noteRange = $notePara.Range
    $noteRange.Text = "<Please refer naming conventions for global or local suffix>"
    $noteRange.Font.Color = 255
}

# ---------------------------------------------------------------------------
# 3) Add an extra blank "ListParagraph" line before "Macros shall NOT be used
#    to replace any keywords of C or brackets." (between it and the existing
#    blank line after "Macros shall all be in upper case").
# ---------------------------------------------------------------------------
$macroRange = $d.Content
$macroFound = $macroRange.Find.Execute(
    "Macros shall all be in upper case",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($macroFound) {
    $macroPara = $macroRange.Paragraphs(1)
    $blankPara = $macroPara.Next()
    $blankPara.Range.InsertParagraphAfter()
}
